# Delete the "Soker et al., 1996 (HUVEC)" row from the VEGFA165_NRP1 sheet.
# That row is row 3 (A3:D3) on the VEGFA165_NRP1 worksheet. Deleting the
# entire row shifts the remaining rows up, removes the now-unused shared
# string from the table, and re-points the other sheets' cells that used
# later shared-string indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("VEGFA165_NRP1")

# Bring this sheet to the front / make it the active tab, matching the
# saved workbook view (activeTab points at VEGFA165_NRP1 afterwards).
$ws.Activate()

# Remove the whole row (shifts rows 4-8 up into rows 3-7).
$ws.Rows(3).Delete()

# Leave the selection on the row that is now in the former row-3's place,
# matching the saved selection state.
$ws.Range("A3:D3").Select()
